$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. "besoin" block: add Status "FIni" for row 6 ---
$ws.Range("E6").Value = "FIni"

# --- 2. "equivalence_besoin" block: add Status "Fini" for row 10 ---
$ws.Range("E10").Value = "Fini"

# --- 3. "don" block: rename fields and drop the duplicated id_ville line ---
$ws.Range("C25").Value = "date"
$ws.Range("C26").Value = "prix_total"
$ws.Range("C27").ClearContents()

# --- 4. insert 3 new rows before "Modeles" (row 33) to make room for the new
#        "dispatch" block; row 32 is already a blank separator row that we
#        reuse as the block header. ---
$ws.Rows("33:35").Insert()

$ws.Range("B32").Value = "dispatch"
$ws.Range("D32").Value = "Christian"
$ws.Range("C33").Value = "id_donnation"
$ws.Range("C34").Value = "id_ville"
$ws.Range("C35").Value = "date_dispatch"

# --- 5. update the view state to match the saved selection/scroll position ---
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Range("E32").Select()
